$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value in column F, row 2 (next to the existing data in B2:E2)
$ws.Range("F2").Value = 2

# Move the active selection to the newly entered cell, matching the author's workflow
$ws.Range("F2").Select()
